# Users_Info.xlsx edit:
#  - Rename sheet "工作表2" -> "Venues"
#  - Populate the (previously empty) Venues sheet with a one-column table
#    of venue codes, starting at row 7
#  - Turn that range into an Excel Table ("表格2" / "Venues" column)
#  - Make the Venues sheet the active/selected sheet, with a specific
#    cell selected on each sheet (mirrors the author's last saved view)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Users")
$ws2 = $wb.Worksheets.Item("工作表2")

# --- rename 工作表2 -> Venues -------------------------------------------
$ws2.Name = "Venues"

# --- fill in the venue codes starting at A7 -----------------------------
$ws2.Range("A7").Value = "Venues"

$venueCodes = @(111,112,113,114,115,116,117,118,119,211,212,213,214,215,216,217,218,219)
for ($i = 0; $i -lt $venueCodes.Length; $i++) {
    $ws2.Cells.Item(8 + $i, 1).Value = $venueCodes[$i]
}

# --- turn A7:A25 into a table --------------------------------------------
$tbl = $ws2.ListObjects.Add(1, $ws2.Range("A7:A25"), 0, 1)
$tbl.Name = "表格2"
$tbl.TableStyle = "TableStyleMedium4"

# --- selections / active sheet -------------------------------------------
[void]$ws1.Range("B1").Select()

$ws2.Activate()
[void]$ws2.Range("F22").Select()
